$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: existing "Pre-Cond" reading becomes "Pre-Cond - 1"; fix the
#     calibration intercept (C2) which had been left as stray text "152.65029."
#     and should really be the numeric value used by the regression.
$ws.Range("C2").Value = 152.65029000000001
$ws.Range("H2").Value = "Pre-Cond - 1"

# --- Row 3: existing "Post-Cond" reading becomes "Post-Cond - 1"; fill in the
#     intercept (C3) and its standard error (E3) that were missing before.
$ws.Range("C3").Value = 95.634039999999999
$ws.Range("E3").Value = 0.00316
$ws.Range("H3").Value = "Post-Cond - 1"

# --- Second round of data collection: two new event rows (4 and 5), each
#     only carrying the FSR id and the Event label for now.
$ws.Range("A4").Value = "FSR_S2"
$ws.Range("H4").Value = "Pre-Cond - 2"

$ws.Range("A5").Value = "FSR_S2"
$ws.Range("H5").Value = "Post-Cond - 2"

# --- Formatting: give the numeric calibration columns (B:G) on rows 2-3 a
#     centered alignment, matching the rest of the table.
$dataRange = $ws.Range("B2:G3")
$dataRange.HorizontalAlignment = -4108

# --- Move the active selection/view the way the author left it.
$ws.Range("D11").Select()

Write-Output "done"
